$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 110.8604276666667
$ws.Range("H2").Value = 332.581283
$ws.Range("I2").Value = 0.2509786052589675
$ws.Range("J2").Value = 0.2509786052589675
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8478306666666667
$ws.Range("N2").Value = 2.543492
$ws.Range("O2").Value = 0.01460351867535248
$ws.Range("P2").Value = 0.01460351867535248
$ws.Range("Q2").Value = 93.99087029558179
$ws.Range("R2").Value = 845.9178326602359
$ws.Range("S2").Value = 0.00366517074901325
$ws.Range("T2").Value = 0.00366517074901325
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 110.8604276666667
$ws.Range("H3").Value = 332.581283
$ws.Range("I3").Value = 0.2509786052589675
$ws.Range("J3").Value = 0.2509786052589675
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.020353
$ws.Range("N3").Value = 21.061059
$ws.Range("O3").Value = 0.1209225617494376
$ws.Range("P3").Value = 0.1209225617494376
$ws.Range("Q3").Value = 778.2793359509664
$ws.Range("R3").Value = 7004.514023558697
$ws.Range("S3").Value = 0.03034897589221521
$ws.Range("T3").Value = 0.03034897589221521
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 110.8604276666667
$ws.Range("H4").Value = 332.581283
$ws.Range("I4").Value = 0.2509786052589675
$ws.Range("J4").Value = 0.2509786052589675
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.462291666666667
$ws.Range("N4").Value = 4.386875
$ws.Range("O4").Value = 0.02518734518879435
$ws.Range("P4").Value = 0.02518734518879435
$ws.Range("Q4").Value = 162.1102795400694
$ws.Range("R4").Value = 1458.992515860625
$ws.Range("S4").Value = 0.00632148476565977
$ws.Range("T4").Value = 0.00632148476565977
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 110.8604276666667
$ws.Range("H5").Value = 332.581283
$ws.Range("I5").Value = 0.2509786052589675
$ws.Range("J5").Value = 0.2509786052589675
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 48.72612633333333
$ws.Range("N5").Value = 146.178379
$ws.Range("O5").Value = 0.8392865743864156
$ws.Range("P5").Value = 0.8392865743864156
$ws.Range("Q5").Value = 5401.799203853362
$ws.Range("R5").Value = 48616.19283468026
$ws.Range("S5").Value = 0.2106429738520793
$ws.Range("T5").Value = 0.2106429738520793
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 184.841802
$ws.Range("H6").Value = 554.525406
$ws.Range("I6").Value = 0.4184661617850055
$ws.Range("J6").Value = 0.4184661617850055
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.8478306666666667
$ws.Range("N6").Value = 2.543492
$ws.Range("O6").Value = 0.01460351867535248
$ws.Range("P6").Value = 0.01460351867535248
$ws.Range("Q6").Value = 156.714548217528
$ws.Range("R6").Value = 1410.430933957752
$ws.Range("S6").Value = 0.006111078408630401
$ws.Range("T6").Value = 0.0061110784086304
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 184.841802
$ws.Range("H7").Value = 554.525406
$ws.Range("I7").Value = 0.4184661617850055
$ws.Range("J7").Value = 0.4184661617850055
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.020353
$ws.Range("N7").Value = 21.061059
$ws.Range("O7").Value = 0.1209225617494376
$ws.Range("P7").Value = 0.1209225617494376
$ws.Range("Q7").Value = 1297.654699196106
$ws.Range("R7").Value = 11678.89229276495
$ws.Range("S7").Value = 0.05060200028849745
$ws.Range("T7").Value = 0.05060200028849745
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 184.841802
$ws.Range("H8").Value = 554.525406
$ws.Range("I8").Value = 0.4184661617850055
$ws.Range("J8").Value = 0.4184661617850055
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.462291666666667
$ws.Range("N8").Value = 4.386875
$ws.Range("O8").Value = 0.02518734518879435
$ws.Range("P8").Value = 0.02518734518879435
$ws.Range("Q8").Value = 270.29262671625
$ws.Range("R8").Value = 2432.63364044625
$ws.Range("S8").Value = 0.01054005166670879
$ws.Range("T8").Value = 0.01054005166670879
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 184.841802
$ws.Range("H9").Value = 554.525406
$ws.Range("I9").Value = 0.4184661617850055
$ws.Range("J9").Value = 0.4184661617850055
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 48.72612633333333
$ws.Range("N9").Value = 146.178379
$ws.Range("O9").Value = 0.8392865743864156
$ws.Range("P9").Value = 0.8392865743864156
$ws.Range("Q9").Value = 9006.624995932985
$ws.Range("R9").Value = 81059.62496339687
$ws.Range("S9").Value = 0.3512130314211688
$ws.Range("T9").Value = 0.3512130314211688
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 95.23175666666667
$ws.Range("H10").Value = 285.69527
$ws.Range("I10").Value = 0.2155966197102082
$ws.Range("J10").Value = 0.2155966197102082
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8478306666666667
$ws.Range("N10").Value = 2.543492
$ws.Range("O10").Value = 0.01460351867535248
$ws.Range("P10").Value = 0.01460351867535248
$ws.Range("Q10").Value = 80.74040374253778
$ws.Range("R10").Value = 726.66363368284
$ws.Range("S10").Value = 0.003148469262280893
$ws.Range("T10").Value = 0.003148469262280893
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 95.23175666666667
$ws.Range("H11").Value = 285.69527
$ws.Range("I11").Value = 0.2155966197102082
$ws.Range("J11").Value = 0.2155966197102082
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.020353
$ws.Range("N11").Value = 21.061059
$ws.Range("O11").Value = 0.1209225617494376
$ws.Range("P11").Value = 0.1209225617494376
$ws.Range("Q11").Value = 668.5605486101034
$ws.Range("R11").Value = 6017.04493749093
$ws.Range("S11").Value = 0.02607049555987766
$ws.Range("T11").Value = 0.02607049555987766
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 95.23175666666667
$ws.Range("H12").Value = 285.69527
$ws.Range("I12").Value = 0.2155966197102082
$ws.Range("J12").Value = 0.2155966197102082
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.462291666666667
$ws.Range("N12").Value = 4.386875
$ws.Range("O12").Value = 0.02518734518879435
$ws.Range("P12").Value = 0.02518734518879435
$ws.Range("Q12").Value = 139.2566041756944
$ws.Range("R12").Value = 1253.30943758125
$ws.Range("S12").Value = 0.005430306482178238
$ws.Range("T12").Value = 0.005430306482178239
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 95.23175666666667
$ws.Range("H13").Value = 285.69527
$ws.Range("I13").Value = 0.2155966197102082
$ws.Range("J13").Value = 0.2155966197102082
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 48.72612633333333
$ws.Range("N13").Value = 146.178379
$ws.Range("O13").Value = 0.8392865743864156
$ws.Range("P13").Value = 0.8392865743864156
$ws.Range("Q13").Value = 4640.274606285259
$ws.Range("R13").Value = 41762.47145656733
$ws.Range("S13").Value = 0.1809473484058714
$ws.Range("T13").Value = 0.1809473484058715
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 50.778675
$ws.Range("H14").Value = 152.336025
$ws.Range("I14").Value = 0.1149586132458188
$ws.Range("J14").Value = 0.1149586132458188
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.8478306666666667
$ws.Range("N14").Value = 2.543492
$ws.Range("O14").Value = 0.01460351867535248
$ws.Range("P14").Value = 0.01460351867535248
$ws.Range("Q14").Value = 43.0517178777
$ws.Range("R14").Value = 387.4654608993
$ws.Range("S14").Value = 0.001678800255427938
$ws.Range("T14").Value = 0.001678800255427938
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 50.778675
$ws.Range("H15").Value = 152.336025
$ws.Range("I15").Value = 0.1149586132458188
$ws.Range("J15").Value = 0.1149586132458188
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 7.020353
$ws.Range("N15").Value = 21.061059
$ws.Range("O15").Value = 0.1209225617494376
$ws.Range("P15").Value = 0.1209225617494376
$ws.Range("Q15").Value = 356.484223372275
$ws.Range("R15").Value = 3208.358010350475
$ws.Range("S15").Value = 0.01390109000884723
$ws.Range("T15").Value = 0.01390109000884724
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 50.778675
$ws.Range("H16").Value = 152.336025
$ws.Range("I16").Value = 0.1149586132458188
$ws.Range("J16").Value = 0.1149586132458188
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.462291666666667
$ws.Range("N16").Value = 4.386875
$ws.Range("O16").Value = 0.02518734518879435
$ws.Range("P16").Value = 0.02518734518879435
$ws.Range("Q16").Value = 74.253233296875
$ws.Range("R16").Value = 668.279099671875
$ws.Range("S16").Value = 0.002895502274247544
$ws.Range("T16").Value = 0.002895502274247544
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 50.778675
$ws.Range("H17").Value = 152.336025
$ws.Range("I17").Value = 0.1149586132458188
$ws.Range("J17").Value = 0.1149586132458188
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 48.72612633333333
$ws.Range("N17").Value = 146.178379
$ws.Range("O17").Value = 0.8392865743864156
$ws.Range("P17").Value = 0.8392865743864156
$ws.Range("Q17").Value = 2474.248133089275
$ws.Range("R17").Value = 22268.23319780348
$ws.Range("S17").Value = 0.09648322070729608
$ws.Range("T17").Value = 0.09648322070729608
